$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 13250
$ws.Range("E2").Value = 1104
$ws.Range("F2").Value = 1104
$ws.Range("G2").Value = 1174
$ws.Range("H2").Value = 867
$ws.Range("I2").Value = 867
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 5440
$ws.Range("L2").Value = 2060
$ws.Range("M2").Value = 3380
$ws.Range("N2").Value = 3379
$ws.Range("O2").Value = 2
$ws.Range("P2").Value = 235
$ws.Range("Q2").Value = 499
$ws.Range("R2").Value = -173
$ws.Range("S2").Value = -214
$ws.Range("T2").Value = 587
$ws.Range("U2").Value = -88
$ws.Range("V2").Value = 189
$ws.Range("W2").Value = 8.33
$ws.Range("X2").Value = 6.54
$ws.Range("Y2").Value = 28.79
$ws.Range("Z2").Value = 16.51
$ws.Range("AA2").Value = 60.93
$ws.Range("AB2").Value = 1551.99
$ws.Range("AC2").Value = 3682
$ws.Range("AD2").Value = 31.09
$ws.Range("AE2").Value = 18688
$ws.Range("AF2").Value = 6.13
$ws.Range("AG2").Value = 850
$ws.Range("AH2").Value = 0.74
$ws.Range("AI2").Value = 17.73
$ws.Range("AJ2").Value = 23533928
$ws.Range("D3").Value = 17105
$ws.Range("E3").Value = 1467
$ws.Range("F3").Value = 1467
$ws.Range("G3").Value = 1553
$ws.Range("H3").Value = 1147
$ws.Range("I3").Value = 1147
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 7688
$ws.Range("L3").Value = 3337
$ws.Range("M3").Value = 4352
$ws.Range("N3").Value = 4351
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 235
$ws.Range("Q3").Value = 2040
$ws.Range("R3").Value = -1591
$ws.Range("S3").Value = -127
$ws.Range("T3").Value = 351
$ws.Range("U3").Value = 1690
$ws.Range("V3").Value = 224
$ws.Range("W3").Value = 8.58
$ws.Range("X3").Value = 6.71
$ws.Range("Y3").Value = 29.67
$ws.Range("Z3").Value = 17.48
$ws.Range("AA3").Value = 76.67
$ws.Range("AB3").Value = 1958.8
$ws.Range("AC3").Value = 4872
$ws.Range("AD3").Value = 47.52
$ws.Range("AE3").Value = 24066
$ws.Range("AF3").Value = 9.619999999999999
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 0.43
$ws.Range("AI3").Value = 15.77
$ws.Range("AJ3").Value = 23533928
$ws.Range("D4").Value = 19345
$ws.Range("E4").Value = 1596
$ws.Range("F4").Value = 1596
$ws.Range("G4").Value = 1701
$ws.Range("H4").Value = 1275
$ws.Range("I4").Value = 1275
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 8994
$ws.Range("L4").Value = 3579
$ws.Range("M4").Value = 5415
$ws.Range("N4").Value = 5414
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 235
$ws.Range("Q4").Value = 1540
$ws.Range("R4").Value = -534
$ws.Range("S4").Value = -201
$ws.Range("T4").Value = 315
$ws.Range("U4").Value = 1225
$ws.Range("V4").Value = 200
$ws.Range("W4").Value = 8.25
$ws.Range("X4").Value = 6.59
$ws.Range("Y4").Value = 26.12
$ws.Range("Z4").Value = 15.29
$ws.Range("AA4").Value = 66.09
$ws.Range("AB4").Value = 2410.56
$ws.Range("AC4").Value = 5419
$ws.Range("AD4").Value = 36.72
$ws.Range("AE4").Value = 29947
$ws.Range("AF4").Value = 6.64
$ws.Range("AG4").Value = 1100
$ws.Range("AH4").Value = 0.55
$ws.Range("AI4").Value = 15.59
$ws.Range("AJ4").Value = 23533928
$ws.Range("D5").Value = 20625
$ws.Range("E5").Value = 1405
$ws.Range("F5").Value = 1405
$ws.Range("G5").Value = 1382
$ws.Range("H5").Value = 959
$ws.Range("I5").Value = 960
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 9757
$ws.Range("L5").Value = 4681
$ws.Range("M5").Value = 5076
$ws.Range("N5").Value = 5075
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 235
$ws.Range("Q5").Value = 662
$ws.Range("R5").Value = -1543
$ws.Range("S5").Value = 172
$ws.Range("T5").Value = 1166
$ws.Range("U5").Value = -504
$ws.Range("V5").Value = 1507
$ws.Range("W5").Value = 6.81
$ws.Range("X5").Value = 4.65
$ws.Range("Y5").Value = 18.3
$ws.Range("Z5").Value = 10.23
$ws.Range("AA5").Value = 92.20999999999999
$ws.Range("AB5").Value = 2718.45
$ws.Range("AC5").Value = 4077
$ws.Range("AD5").Value = 44.15
$ws.Range("AE5").Value = 28788
$ws.Range("AF5").Value = 6.25
$ws.Range("AG5").Value = 1200
$ws.Range("AH5").Value = 0.67
$ws.Range("AI5").Value = 22.05
$ws.Range("AJ5").Value = 23533928
$ws.Range("D6").Value = 19285
$ws.Range("E6").Value = 560
$ws.Range("F6").Value = 560
$ws.Range("G6").Value = 1247
$ws.Range("H6").Value = 900
$ws.Range("I6").Value = 900
$ws.Range("K6").Value = 9402
$ws.Range("L6").Value = 3657
$ws.Range("M6").Value = 5745
$ws.Range("N6").Value = 5744
$ws.Range("P6").Value = 235
$ws.Range("Q6").Value = 637
$ws.Range("R6").Value = 501
$ws.Range("S6").Value = -1079
$ws.Range("T6").Value = 119
$ws.Range("U6").Value = 518
$ws.Range("V6").Value = 639
$ws.Range("W6").Value = 2.91
$ws.Range("X6").Value = 4.67
$ws.Range("Y6").Value = 16.63
$ws.Range("Z6").Value = 9.390000000000001
$ws.Range("AA6").Value = 63.66
$ws.Range("AB6").Value = 2999.76
$ws.Range("AC6").Value = 3823
$ws.Range("AD6").Value = 17.63
$ws.Range("AE6").Value = 32581
$ws.Range("AF6").Value = 2.07
$ws.Range("AG6").Value = 1200
$ws.Range("AH6").Value = 1.78
$ws.Range("AI6").Value = 23.51
$ws.Range("AJ6").Value = 23533928
$ws.Range("D7").Value = 17139
$ws.Range("E7").Value = 521
$ws.Range("G7").Value = 656
$ws.Range("H7").Value = 355
$ws.Range("I7").Value = 382
$ws.Range("K7").Value = 10253
$ws.Range("L7").Value = 4350
$ws.Range("M7").Value = 5903
$ws.Range("N7").Value = 5905
$ws.Range("P7").Value = 238
$ws.Range("Q7").Value = 712
$ws.Range("R7").Value = -285
$ws.Range("S7").Value = -276
$ws.Range("T7").Value = 183
$ws.Range("U7").Value = 686
$ws.Range("W7").Value = 3.04
$ws.Range("X7").Value = 2.07
$ws.Range("Y7").Value = 6.56
$ws.Range("Z7").Value = 3.62
$ws.Range("AA7").Value = 73.68000000000001
$ws.Range("AC7").Value = 1624
$ws.Range("AD7").Value = 43.89
$ws.Range("AE7").Value = 33494
$ws.Range("AF7").Value = 2.13
$ws.Range("AG7").Value = 1195
$ws.Range("AH7").Value = 1.68
$ws.Range("AI7").Value = 73.56999999999999
$ws.Range("D8").Value = 18298
$ws.Range("E8").Value = 834
$ws.Range("G8").Value = 948
$ws.Range("H8").Value = 665
$ws.Range("I8").Value = 680
$ws.Range("K8").Value = 10710
$ws.Range("L8").Value = 4384
$ws.Range("M8").Value = 6326
$ws.Range("N8").Value = 6319
$ws.Range("P8").Value = 238
$ws.Range("Q8").Value = 941
$ws.Range("R8").Value = -408
$ws.Range("S8").Value = -347
$ws.Range("T8").Value = 277
$ws.Range("U8").Value = 601
$ws.Range("W8").Value = 4.56
$ws.Range("X8").Value = 3.64
$ws.Range("Y8").Value = 11.13
$ws.Range("Z8").Value = 6.35
$ws.Range("AA8").Value = 69.31
$ws.Range("AC8").Value = 2890
$ws.Range("AD8").Value = 24.67
$ws.Range("AE8").Value = 35845
$ws.Range("AF8").Value = 1.99
$ws.Range("AG8").Value = 1200
$ws.Range("AH8").Value = 1.68
$ws.Range("AI8").Value = 41.52
$ws.Range("D9").Value = 19772
$ws.Range("E9").Value = 1082
$ws.Range("G9").Value = 1215
$ws.Range("H9").Value = 870
$ws.Range("I9").Value = 871
$ws.Range("K9").Value = 11396
$ws.Range("L9").Value = 4419
$ws.Range("M9").Value = 6977
$ws.Range("N9").Value = 6938
$ws.Range("P9").Value = 238
$ws.Range("Q9").Value = 1117
$ws.Range("R9").Value = -413
$ws.Range("S9").Value = -275
$ws.Range("T9").Value = 335
$ws.Range("U9").Value = 758
$ws.Range("W9").Value = 5.47
$ws.Range("X9").Value = 4.4
$ws.Range("Y9").Value = 13.14
$ws.Range("Z9").Value = 7.87
$ws.Range("AA9").Value = 63.34
$ws.Range("AC9").Value = 3702
$ws.Range("AD9").Value = 19.26
$ws.Range("AE9").Value = 39353
$ws.Range("AF9").Value = 1.81
$ws.Range("AG9").Value = 1206
$ws.Range("AH9").Value = 1.69
$ws.Range("AI9").Value = 32.57
